$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (new quarter columns), shifting
# old D:K data to F:M
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number formats / styles from column F (the old column D, now
# shifted) into the two new columns D:E so dates/numbers render the same
$ws.Range("F5:F102").Copy() | Out-Null
$ws.Range("D5:E102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Rows 5, 6, 37 and 79 are single-cell section headers (columns A/B only)
# and should stay that way - undo the formatting paste on those rows
$ws.Range("D5:E6").Clear()
$ws.Range("D37:E37").Clear()
$ws.Range("D79:E79").Clear()

# Make the two new columns the same width as the rest of the data columns
$ws.Range("D1:M1").EntireColumn.ColumnWidth = 15.1666666666667

# Populate the two newly inserted columns (D = latest quarter, E = prior
# quarter) with the restated / new financial figures
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 17438000
$ws.Range("E8").Value = 16690000
$ws.Range("D9").Value = 8776000
$ws.Range("E9").Value = 8420000
$ws.Range("D10").Value = 8662000
$ws.Range("E10").Value = 8270000
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 177000
$ws.Range("E14").Value = 137000
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 13551000
$ws.Range("E17").Value = 13137000
$ws.Range("D18").Value = 3887000
$ws.Range("E18").Value = 3553000
$ws.Range("D20").Value = 167000
$ws.Range("E20").Value = 516000
$ws.Range("D21").Value = 4704000
$ws.Range("E21").Value = 4712000
$ws.Range("D22").Value = 138000
$ws.Range("E22").Value = 129000
$ws.Range("D23").Value = 3916000
$ws.Range("E23").Value = 3940000
$ws.Range("D24").Value = 700000
$ws.Range("E24").Value = 729000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 3216000
$ws.Range("E26").Value = 3211000
$ws.Range("D27").Value = 3129000
$ws.Range("E27").Value = 3133000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -167000
$ws.Range("E32").Value = -516000
$ws.Range("D33").Value = 3129000
$ws.Range("E33").Value = 3133000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 3129000
$ws.Range("E35").Value = 3133000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 3696000
$ws.Range("E41").Value = 2545000
$ws.Range("D42").Value = 8421000
$ws.Range("E42").Value = 8708000
$ws.Range("D43").Value = 5055000
$ws.Range("E43").Value = 5035000
$ws.Range("D44").Value = 5281000
$ws.Range("E44").Value = 5182000
$ws.Range("D45").Value = 1978000
$ws.Range("E45").Value = 1876000
$ws.Range("D46").Value = 24431000
$ws.Range("E46").Value = 23346000
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 20822000
$ws.Range("E48").Value = 20590000
$ws.Range("D49").Value = 72879000
$ws.Range("E49").Value = 69144000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 5555000
$ws.Range("E52").Value = 5360000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 123687000
$ws.Range("E54").Value = 118440000
$ws.Range("D57").Value = 10266000
$ws.Range("E57").Value = 10243000
$ws.Range("D58").Value = 12113000
$ws.Range("E58").Value = 10508000
$ws.Range("D59").Value = 8868000
$ws.Range("E59").Value = 8469000
$ws.Range("D60").Value = 31247000
$ws.Range("E60").Value = 29220000
$ws.Range("D61").Value = 21514000
$ws.Range("E61").Value = 20779000
$ws.Range("D62").Value = 16483000
$ws.Range("E62").Value = 15937000
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 69697000
$ws.Range("E66").Value = 66204000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 946000
$ws.Range("E70").Value = 951000
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 101170000
$ws.Range("E72").Value = 99831000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 53044000
$ws.Range("E76").Value = 51285000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 3129000
$ws.Range("E81").Value = 3133000
$ws.Range("D83").Value = 650000
$ws.Range("E83").Value = 643000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 4007000
$ws.Range("E89").Value = 3567000
$ws.Range("D91").Value = -701000
$ws.Range("E91").Value = -1080000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -3845000
$ws.Range("E94").Value = -865000
$ws.Range("D96").Value = -1850000
$ws.Range("E96").Value = -1853000
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 1032000
$ws.Range("E100").Value = -2656000
$ws.Range("D101").Value = -43000
$ws.Range("E101").Value = -70000
$ws.Range("D102").Value = 1151000
$ws.Range("E102").Value = -24000

# A handful of historical figures were also restated in this edit
$ws.Range("H9").Value = 17270000
$ws.Range("I9").Value = 16380000
$ws.Range("H10").Value = 125000
$ws.Range("I10").Value = 273000
$ws.Range("H17").Value = 13476000
$ws.Range("I17").Value = 13005000
$ws.Range("H18").Value = 3919000
$ws.Range("I18").Value = 3648000
$ws.Range("H20").Value = 358000
$ws.Range("I20").Value = 218000
$ws.Range("H21").Value = 4953000
$ws.Range("H22").Value = 244000
$ws.Range("H32").Value = -358000
$ws.Range("I32").Value = -218000
